$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 5.6
$ws.Range("J2").Value = 4.4
$ws.Range("P2").Value = 2.46
$ws.Range("S2").Value = 2.64
$ws.Range("AK2").Value = 16
$ws.Range("G3").Value = 1.46
$ws.Range("H3").Value = 9.6
$ws.Range("I3").Value = 11.5
$ws.Range("N3").Value = 3.5
$ws.Range("Q3").Value = 2.02
$ws.Range("U3").Value = 1.68
$ws.Range("V3").Value = 1.09
$ws.Range("W3").Value = 3.15
$ws.Range("X3").Value = 17
$ws.Range("G4").Value = 4.4
$ws.Range("H4").Value = 1.94
$ws.Range("I4").Value = 1.95
$ws.Range("K4").Value = 4
$ws.Range("S4").Value = 2.92
$ws.Range("U4").Value = 2.32
$ws.Range("V4").Value = 2.04
$ws.Range("W4").Value = 1.3
$ws.Range("AB4").Value = 19
$ws.Range("AI4").Value = 30
$ws.Range("AL4").Value = 50
$ws.Range("AO4").Value = 11
$ws.Range("F5").Value = 1.65
$ws.Range("G5").Value = 1.66
$ws.Range("H5").Value = 5.6
$ws.Range("I5").Value = 5.8
$ws.Range("K5").Value = 4.6
$ws.Range("F6").Value = 1.42
$ws.Range("J6").Value = 5.5
$ws.Range("K6").Value = 6.4
$ws.Range("Q6").Value = 1.31
$ws.Range("S6").Value = 1.79
$ws.Range("AA6").Value = 190
$ws.Range("AG6").Value = 12.5
$ws.Range("AN6").Value = 4
$ws.Range("F7").Value = 1.98
$ws.Range("G7").Value = 2.14
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 3.7
$ws.Range("J7").Value = 4.1
$ws.Range("K7").Value = 4.7
$ws.Range("O7").Value = 1.16
$ws.Range("P7").Value = 2.72
$ws.Range("Q7").Value = 1.48
$ws.Range("T7").Value = 1.49
$ws.Range("U7").Value = 2.62
$ws.Range("W7").Value = 1.87
$ws.Range("AA7").Value = 65
$ws.Range("AC7").Value = 13
$ws.Range("AH7").Value = 980
$ws.Range("H8").Value = 3.5
$ws.Range("N8").Value = 5.4
$ws.Range("P8").Value = 2.5
$ws.Range("R8").Value = 1.6
$ws.Range("S8").Value = 2.4
$ws.Range("U8").Value = 2.52
$ws.Range("AD8").Value = 18
$ws.Range("AE8").Value = 42
$ws.Range("AN8").Value = 10.5
$ws.Range("G9").Value = 1.51
$ws.Range("H9").Value = 6.6
$ws.Range("I9").Value = 7.8
$ws.Range("J9").Value = 5.2
$ws.Range("K9").Value = 5.9
$ws.Range("L9").Value = 1.24
$ws.Range("P9").Value = 2.96
$ws.Range("R9").Value = 1.79
$ws.Range("W9").Value = 2.96
$ws.Range("AO9").Value = 1000
$ws.Range("H10").Value = 4
$ws.Range("N10").Value = 2.62
$ws.Range("U10").Value = 1.76
$ws.Range("AJ10").Value = 980
$ws.Range("AK10").Value = 36
$ws.Range("AL10").Value = 60
$ws.Range("AM10").Value = 230
$ws.Range("G11").Value = 2.98
$ws.Range("I11").Value = 2.82
$ws.Range("K11").Value = 3.7
$ws.Range("O11").Value = 1.31
$ws.Range("P11").Value = 1.96
$ws.Range("R11").Value = 1.37
$ws.Range("S11").Value = 3.3
$ws.Range("T11").Value = 1.72
$ws.Range("U11").Value = 2.2
$ws.Range("W11").Value = 1.5
$ws.Range("X11").Value = 18
$ws.Range("AB11").Value = 14.5
$ws.Range("AN11").Value = 32
$ws.Range("F12").Value = 1.79
$ws.Range("I12").Value = 6.8
$ws.Range("K12").Value = 3.8
$ws.Range("S12").Value = 3.5
$ws.Range("T12").Value = 2.04
$ws.Range("F13").Value = 1.92
$ws.Range("G13").Value = 1.94
$ws.Range("I13").Value = 4.5
$ws.Range("R13").Value = 1.5
$ws.Range("T13").Value = 1.71
$ws.Range("W13").Value = 2.06
$ws.Range("AD13").Value = 16.5
$ws.Range("AH13").Value = 16.5
$ws.Range("AO13").Value = 44
$ws.Range("F14").Value = 3.35
$ws.Range("G14").Value = 3.45
$ws.Range("H14").Value = 2.18
$ws.Range("I14").Value = 2.22
$ws.Range("P14").Value = 2.48
$ws.Range("Q14").Value = 1.63
$ws.Range("R14").Value = 1.6
$ws.Range("S14").Value = 2.58
$ws.Range("V14").Value = 1.82
$ws.Range("W14").Value = 1.41
$ws.Range("Y14").Value = 14
$ws.Range("AA14").Value = 27
$ws.Range("AB14").Value = 18
$ws.Range("AF14").Value = 27
$ws.Range("AO14").Value = 11.5
$ws.Range("H15").Value = 2.8
$ws.Range("I15").Value = 2.84
$ws.Range("R15").Value = 1.89
$ws.Range("V15").Value = 1.54
$ws.Range("AO15").Value = 12
$ws.Range("T16").Value = 1.87
$ws.Range("Y16").Value = 65
$ws.Range("AA16").Value = 580
$ws.Range("AI16").Value = 140
$ws.Range("G17").Value = 1.33
$ws.Range("R17").Value = 1.62
$ws.Range("AA17").Value = 460
$ws.Range("AB17").Value = 9.800000000000001
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 3.05
$ws.Range("H18").Value = 2.56
$ws.Range("I18").Value = 2.58
$ws.Range("P18").Value = 2.12
$ws.Range("Q18").Value = 1.86
$ws.Range("V18").Value = 1.63
$ws.Range("W18").Value = 1.48
$ws.Range("AF18").Value = 20
$ws.Range("AN18").Value = 25
$ws.Range("AO18").Value = 19
$ws.Range("O19").Value = 1.37
$ws.Range("Y19").Value = 10
$ws.Range("G20").Value = 5
$ws.Range("G21").Value = 3.25
$ws.Range("H21").Value = 2.58
$ws.Range("O21").Value = 1.41
$ws.Range("V21").Value = 1.56
$ws.Range("W21").Value = 1.45
$ws.Range("AB21").Value = 11
$ws.Range("AF21").Value = 20
$ws.Range("AO21").Value = 32
$ws.Range("R22").Value = 1.22
